$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 35989.355
$ws.Range("I12").Value2 = 300.08334
$ws.Range("J12").Value2 = 250125
$ws.Range("K12").Value2 = 300.08334
$ws.Range("L12").Value2 = 250125
$ws.Range("M12").Value2 = -130.08334
$ws.Range("N12").Value2 = -250465

$ws.Range("H49").Value2 = 495
$ws.Range("I49").Value2 = 0
$ws.Range("J49").Value2 = 495
$ws.Range("K49").Value2 = 0
$ws.Range("L49").Value2 = 1485
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value2 = -1757

$ws.Range("H55").Value2 = 239.8
$ws.Range("I55").Value2 = 71
$ws.Range("K55").Value2 = 71
$ws.Range("M55").Value2 = 143

$ws.Range("H106").Value2 = 1578.3334
$ws.Range("I106").Value2 = 1675.625
$ws.Range("J106").Value2 = 800
$ws.Range("K106").Value2 = 1675.625
$ws.Range("L106").Value2 = 800
$ws.Range("M106").Value2 = -1044.625
$ws.Range("N106").Value2 = -2062

$ws.Range("H111").Value2 = 5563880
$ws.Range("I111").Value2 = 11484.917
$ws.Range("J111").Value2 = 16668671
$ws.Range("K111").Value2 = 34454.751
$ws.Range("L111").Value2 = 50006013
$ws.Range("M111").Value2 = -31387.751
$ws.Range("N111").Value2 = -50012147

$ws.Range("H112").Value2 = 1081.9333
$ws.Range("J112").Value2 = 1172.2307
$ws.Range("L112").Value2 = 3516.6921
$ws.Range("N112").Value2 = -5732.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value2 = 3787.8
$ws.Range("I122").Value2 = 2993.3333
$ws.Range("J122").Value2 = 4979.5
$ws.Range("K122").Value2 = 8979.999899999999
$ws.Range("L122").Value2 = 14938.5
$ws.Range("M122").Value2 = -6529.999899999999
$ws.Range("N122").Value2 = -19838.5

$ws.Range("H132").Value2 = 1958.4906
$ws.Range("I132").Value2 = 1965.6578
$ws.Range("J132").Value2 = 1940.3334
$ws.Range("K132").Value2 = 5896.9734
$ws.Range("L132").Value2 = 5821.0002
$ws.Range("M132").Value2 = -3366.9734
$ws.Range("N132").Value2 = -10881.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value2 = 12611.875
$ws.Range("I19").Value2 = 31.666666
$ws.Range("J19").Value2 = 20160
$ws.Range("K19").Value2 = 31.666666
$ws.Range("L19").Value2 = 20160
$ws.Range("M19").Value2 = 138.333334
$ws.Range("N19").Value2 = -20500

$ws.Range("H22").Value2 = 395.85715
$ws.Range("I22").Value2 = 354.2
$ws.Range("J22").Value2 = 500
$ws.Range("K22").Value2 = 354.2
$ws.Range("L22").Value2 = 500
$ws.Range("M22").Value2 = -4.199999999999989
$ws.Range("N22").Value2 = -1200

$ws.Range("H24").Value2 = 12611.875
$ws.Range("I24").Value2 = 31.666666
$ws.Range("J24").Value2 = 20160
$ws.Range("K24").Value2 = 31.666666
$ws.Range("L24").Value2 = 20160
$ws.Range("M24").Value2 = 138.333334
$ws.Range("N24").Value2 = -20500

$ws.Range("H29").Value2 = 9750
$ws.Range("J29").Value2 = 9750
$ws.Range("L29").Value2 = 9750
$ws.Range("N29").Value2 = -10336

$ws.Range("H31").Value2 = 2469.0925
$ws.Range("I31").Value2 = 1635.25
$ws.Range("J31").Value2 = 2820.1843
$ws.Range("K31").Value2 = 1635.25
$ws.Range("L31").Value2 = 2820.1843
$ws.Range("M31").Value2 = -1340.25
$ws.Range("N31").Value2 = -3410.1843

$ws.Range("H34").Value2 = 2469.0925
$ws.Range("I34").Value2 = 1635.25
$ws.Range("J34").Value2 = 2820.1843
$ws.Range("K34").Value2 = 1635.25
$ws.Range("L34").Value2 = 2820.1843
$ws.Range("M34").Value2 = -1433.25
$ws.Range("N34").Value2 = -3224.1843

$ws.Range("H41").Value2 = 9696.25
$ws.Range("I41").Value2 = 3036.6667
$ws.Range("J41").Value2 = 13692
$ws.Range("K41").Value2 = 3036.6667
$ws.Range("L41").Value2 = 13692
$ws.Range("M41").Value2 = -2608.6667
$ws.Range("N41").Value2 = -14548

$ws.Range("H42").Value2 = 0
$ws.Range("J42").Value2 = 0
$ws.Range("L42").Value2 = 0
$ws.Range("N42").ClearContents()

$ws.Range("H50").Value2 = 10760
$ws.Range("J50").Value2 = 13640
$ws.Range("L50").Value2 = 13640
$ws.Range("N50").Value2 = -14890

$ws.Range("H51").Value2 = 7922.923
$ws.Range("I51").Value2 = 0
$ws.Range("J51").Value2 = 7922.923
$ws.Range("K51").Value2 = 0
$ws.Range("L51").Value2 = 7922.923
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value2 = -9394.922999999999

$ws.Range("H59").Value2 = 26890
$ws.Range("J59").Value2 = 29780
$ws.Range("L59").Value2 = 29780
$ws.Range("N59").Value2 = -32070

$ws.Range("H60").Value2 = 15232.8
$ws.Range("J60").Value2 = 18853.334
$ws.Range("L60").Value2 = 18853.334
$ws.Range("N60").Value2 = -19875.334

$ws.Range("H61").Value2 = 7922.923
$ws.Range("I61").Value2 = 0
$ws.Range("J61").Value2 = 7922.923
$ws.Range("K61").Value2 = 0
$ws.Range("L61").Value2 = 7922.923
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value2 = -8618.922999999999

$ws.Range("H105").Value2 = 1270.8235
$ws.Range("I105").Value2 = 1262.4
$ws.Range("J105").Value2 = 1282.8572
$ws.Range("K105").Value2 = 1262.4
$ws.Range("L105").Value2 = 1282.8572
$ws.Range("M105").Value2 = 484.5999999999999
$ws.Range("N105").Value2 = -4776.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value2 = 499.625
$ws.Range("J34").Value2 = 499.625
$ws.Range("L34").Value2 = 1498.875
$ws.Range("N34").Value2 = -1666.875

$ws.Range("H58").Value2 = 2325
$ws.Range("I58").Value2 = 2316.6667
$ws.Range("J58").Value2 = 2350
$ws.Range("K58").Value2 = 6950.000100000001
$ws.Range("L58").Value2 = 7050
$ws.Range("M58").Value2 = -6822.000100000001
$ws.Range("N58").Value2 = -7306

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value2 = 652
$ws.Range("I3").Value2 = 200
$ws.Range("J3").Value2 = 1104
$ws.Range("K3").Value2 = 200
$ws.Range("L3").Value2 = 1104
$ws.Range("M3").Value2 = -84
$ws.Range("N3").Value2 = -1336

$ws.Range("H5").Value2 = 1669666.6
$ws.Range("J5").Value2 = 4500
$ws.Range("L5").Value2 = 4500
$ws.Range("N5").Value2 = -4724

$ws.Range("H12").Value2 = 6191524
$ws.Range("I12").Value2 = 6842842
$ws.Range("J12").Value2 = 4000
$ws.Range("K12").Value2 = 6842842
$ws.Range("L12").Value2 = 4000
$ws.Range("M12").Value2 = -6842702
$ws.Range("N12").Value2 = -4280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 1446417.4
$ws.Range("I46").Value2 = 404
$ws.Range("J46").Value2 = 2530927.5
$ws.Range("K46").Value2 = 404
$ws.Range("L46").Value2 = 2530927.5
$ws.Range("M46").Value2 = -216
$ws.Range("N46").Value2 = -2531303.5

$ws.Range("H55").Value2 = 438073.94
$ws.Range("I55").Value2 = 948348.8
$ws.Range("J55").Value2 = 695.4286
$ws.Range("K55").Value2 = 948348.8
$ws.Range("L55").Value2 = 695.4286
$ws.Range("M55").Value2 = -948175.8
$ws.Range("N55").Value2 = -1041.4286

$ws.Range("H58").Value2 = 3700
$ws.Range("J58").Value2 = 4950
$ws.Range("L58").Value2 = 4950
$ws.Range("N58").Value2 = -5470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value2 = 60914
$ws.Range("I33").Value2 = 0
$ws.Range("K33").Value2 = 0
$ws.Range("M33").ClearContents()

$ws.Range("H36").Value2 = 60914
$ws.Range("I36").Value2 = 0
$ws.Range("K36").Value2 = 0
$ws.Range("M36").ClearContents()

$ws.Range("H56").Value2 = 33484.75
$ws.Range("I56").Value2 = 4000
$ws.Range("J56").Value2 = 43313
$ws.Range("K56").Value2 = 4000
$ws.Range("L56").Value2 = 43313
$ws.Range("M56").Value2 = -3286
$ws.Range("N56").Value2 = -44741

$ws.Range("H122").Value2 = 2402.2354
$ws.Range("J122").Value2 = 3553.5715
$ws.Range("L122").Value2 = 10660.7145
$ws.Range("N122").Value2 = -15560.7145

$ws.Range("H126").Value2 = 1614.4546
$ws.Range("I126").Value2 = 1375.9
$ws.Range("J126").Value2 = 4000
$ws.Range("K126").Value2 = 4127.700000000001
$ws.Range("L126").Value2 = 12000
$ws.Range("M126").Value2 = -1657.700000000001
$ws.Range("N126").Value2 = -16940
